# Auto-generated edit script: updates NATMI LR-pair computed metrics
# following recomputation with 3 ligand/receptor-expressing cells (was 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 37.05583833333333
$ws.Range("H2").Value = 111.167515
$ws.Range("I2").Value = 0.008431126118266585
$ws.Range("J2").Value = 0.008431126118266585
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.744736333333333
$ws.Range("N2").Value = 20.234209
$ws.Range("O2").Value = 0.01049273071342527
$ws.Range("P2").Value = 0.01049273071342527
$ws.Range("Q2").Value = 249.9318591689594
$ws.Range("R2").Value = 2249.386732520635
$ws.Range("S2").Value = 0.00008846553596989779
$ws.Range("T2").Value = 0.0000884655359698978

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 37.05583833333333
$ws.Range("H3").Value = 111.167515
$ws.Range("I3").Value = 0.008431126118266585
$ws.Range("J3").Value = 0.008431126118266585
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 182.5316823333334
$ws.Range("N3").Value = 547.595047
$ws.Range("O3").Value = 0.283963033503136
$ws.Range("P3").Value = 0.2839630335031361
$ws.Range("Q3").Value = 6763.864511255357
$ws.Range("R3").Value = 60874.7806012982
$ws.Range("S3").Value = 0.002394128148390499
$ws.Range("T3").Value = 0.0023941281483905

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 37.05583833333333
$ws.Range("H4").Value = 111.167515
$ws.Range("I4").Value = 0.008431126118266585
$ws.Range("J4").Value = 0.008431126118266585
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 249.1329143333333
$ws.Range("N4").Value = 747.398743
$ws.Range("O4").Value = 0.3875740210972192
$ws.Range("P4").Value = 0.3875740210972192
$ws.Range("Q4").Value = 9231.828997048184
$ws.Range("R4").Value = 83086.46097343364
$ws.Range("S4").Value = 0.003267685452034369
$ws.Range("T4").Value = 0.003267685452034369

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.05583833333333
$ws.Range("H5").Value = 111.167515
$ws.Range("I5").Value = 0.008431126118266585
$ws.Range("J5").Value = 0.008431126118266585
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 186.9310863333334
$ws.Range("N5").Value = 560.793259
$ws.Range("O5").Value = 0.2908071500393791
$ws.Range("P5").Value = 0.2908071500393791
$ws.Range("Q5").Value = 6926.888114642377
$ws.Range("R5").Value = 62341.99303178139
$ws.Range("S5").Value = 0.002451831758075678
$ws.Range("T5").Value = 0.002451831758075678

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.05583833333333
$ws.Range("H6").Value = 111.167515
$ws.Range("I6").Value = 0.008431126118266585
$ws.Range("J6").Value = 0.008431126118266585
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 17.46044133333334
$ws.Range("N6").Value = 52.38132400000001
$ws.Range("O6").Value = 0.02716306464684043
$ws.Range("P6").Value = 0.02716306464684043
$ws.Range("Q6").Value = 647.0112912766511
$ws.Range("R6").Value = 5823.10162148986
$ws.Range("S6").Value = 0.00022901522379614
$ws.Range("T6").Value = 0.0002290152237961401

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 50.89916233333333
$ws.Range("H7").Value = 152.697487
$ws.Range("I7").Value = 0.01158082710438721
$ws.Range("J7").Value = 0.01158082710438721
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.744736333333333
$ws.Range("N7").Value = 20.234209
$ws.Range("O7").Value = 0.01049273071342527
$ws.Range("P7").Value = 0.01049273071342527
$ws.Range("Q7").Value = 343.3014295258648
$ws.Range("R7").Value = 3089.712865732783
$ws.Range("S7").Value = 0.0001215145002450716
$ws.Range("T7").Value = 0.0001215145002450716

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.89916233333333
$ws.Range("H8").Value = 152.697487
$ws.Range("I8").Value = 0.01158082710438721
$ws.Range("J8").Value = 0.01158082710438721
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 182.5316823333334
$ws.Range("N8").Value = 547.595047
$ws.Range("O8").Value = 0.283963033503136
$ws.Range("P8").Value = 0.2839630335031361
$ws.Range("Q8").Value = 9290.709730060766
$ws.Range("R8").Value = 83616.3875705469
$ws.Range("S8").Value = 0.003288526795037132
$ws.Range("T8").Value = 0.003288526795037133

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.89916233333333
$ws.Range("H9").Value = 152.697487
$ws.Range("I9").Value = 0.01158082710438721
$ws.Range("J9").Value = 0.01158082710438721
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 249.1329143333333
$ws.Range("N9").Value = 747.398743
$ws.Range("O9").Value = 0.3875740210972192
$ws.Range("P9").Value = 0.3875740210972192
$ws.Range("Q9").Value = 12680.65664922876
$ws.Range("R9").Value = 114125.9098430588
$ws.Range("S9").Value = 0.004488427728479018
$ws.Range("T9").Value = 0.004488427728479018

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.89916233333333
$ws.Range("H10").Value = 152.697487
$ws.Range("I10").Value = 0.01158082710438721
$ws.Range("J10").Value = 0.01158082710438721
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 186.9310863333334
$ws.Range("N10").Value = 560.793259
$ws.Range("O10").Value = 0.2908071500393791
$ws.Range("P10").Value = 0.2908071500393791
$ws.Range("Q10").Value = 9514.635708426682
$ws.Range("R10").Value = 85631.72137584013
$ws.Range("S10").Value = 0.003367787325325641
$ws.Range("T10").Value = 0.003367787325325641

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 50.89916233333333
$ws.Range("H11").Value = 152.697487
$ws.Range("I11").Value = 0.01158082710438721
$ws.Range("J11").Value = 0.01158082710438721
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 17.46044133333334
$ws.Range("N11").Value = 52.38132400000001
$ws.Range("O11").Value = 0.02716306464684043
$ws.Range("P11").Value = 0.02716306464684043
$ws.Range("Q11").Value = 888.7218378369765
$ws.Range("R11").Value = 7998.496540532788
$ws.Range("S11").Value = 0.0003145707553003517
$ws.Range("T11").Value = 0.0003145707553003518

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2007.446289
$ws.Range("H12").Value = 6022.338867
$ws.Range("I12").Value = 0.4567440273772037
$ws.Range("J12").Value = 0.4567440273772037
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.744736333333333
$ws.Range("N12").Value = 20.234209
$ws.Range("O12").Value = 0.01049273071342527
$ws.Range("P12").Value = 0.01049273071342527
$ws.Range("Q12").Value = 13539.69592263347
$ws.Range("R12").Value = 121857.2633037012
$ws.Range("S12").Value = 0.004792492084234338
$ws.Range("T12").Value = 0.004792492084234339

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2007.446289
$ws.Range("H13").Value = 6022.338867
$ws.Range("I13").Value = 0.4567440273772037
$ws.Range("J13").Value = 0.4567440273772037
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 182.5316823333334
$ws.Range("N13").Value = 547.595047
$ws.Range("O13").Value = 0.283963033503136
$ws.Range("P13").Value = 0.2839630335031361
$ws.Range("Q13").Value = 366422.5483249769
$ws.Range("R13").Value = 3297802.934924792
$ws.Range("S13").Value = 0.1296984195484701
$ws.Range("T13").Value = 0.1296984195484702

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2007.446289
$ws.Range("H14").Value = 6022.338867
$ws.Range("I14").Value = 0.4567440273772037
$ws.Range("J14").Value = 0.4567440273772037
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 249.1329143333333
$ws.Range("N14").Value = 747.398743
$ws.Range("O14").Value = 0.3875740210972192
$ws.Range("P14").Value = 0.3875740210972192
$ws.Range("Q14").Value = 500120.944346205
$ws.Range("R14").Value = 4501088.499115844
$ws.Range("S14").Value = 0.1770221193027212
$ws.Range("T14").Value = 0.1770221193027212

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2007.446289
$ws.Range("H15").Value = 6022.338867
$ws.Range("I15").Value = 0.4567440273772037
$ws.Range("J15").Value = 0.4567440273772037
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 186.9310863333334
$ws.Range("N15").Value = 560.793259
$ws.Range("O15").Value = 0.2908071500393791
$ws.Range("P15").Value = 0.2908071500393791
$ws.Range("Q15").Value = 375254.1155585887
$ws.Range("R15").Value = 3377287.040027298
$ws.Range("S15").Value = 0.1328244288990727
$ws.Range("T15").Value = 0.1328244288990727

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2007.446289
$ws.Range("H16").Value = 6022.338867
$ws.Range("I16").Value = 0.4567440273772037
$ws.Range("J16").Value = 0.4567440273772037
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 17.46044133333334
$ws.Range("N16").Value = 52.38132400000001
$ws.Range("O16").Value = 0.02716306464684043
$ws.Range("P16").Value = 0.02716306464684043
$ws.Range("Q16").Value = 35050.89815890222
$ws.Range("R16").Value = 315458.08343012
$ws.Range("S16").Value = 0.01240656754270524
$ws.Range("T16").Value = 0.01240656754270524

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2293.273345666667
$ws.Range("H17").Value = 6879.820037
$ws.Range("I17").Value = 0.5217768014597114
$ws.Range("J17").Value = 0.5217768014597114
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.744736333333333
$ws.Range("N17").Value = 20.234209
$ws.Range("O17").Value = 0.01049273071342527
$ws.Range("P17").Value = 0.01049273071342527
$ws.Range("Q17").Value = 15467.52405678286
$ws.Range("R17").Value = 139207.7165110457
$ws.Range("S17").Value = 0.005474863470229114
$ws.Range("T17").Value = 0.005474863470229115

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 2293.273345666667
$ws.Range("H18").Value = 6879.820037
$ws.Range("I18").Value = 0.5217768014597114
$ws.Range("J18").Value = 0.5217768014597114
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 182.5316823333334
$ws.Range("N18").Value = 547.595047
$ws.Range("O18").Value = 0.283963033503136
$ws.Range("P18").Value = 0.2839630335031361
$ws.Range("Q18").Value = 418595.0418347286
$ws.Range("R18").Value = 3767355.376512557
$ws.Range("S18").Value = 0.1481653233540632
$ws.Range("T18").Value = 0.1481653233540632

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 2293.273345666667
$ws.Range("H19").Value = 6879.820037
$ws.Range("I19").Value = 0.5217768014597114
$ws.Range("J19").Value = 0.5217768014597114
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 249.1329143333333
$ws.Range("N19").Value = 747.398743
$ws.Range("O19").Value = 0.3875740210972192
$ws.Range("P19").Value = 0.3875740210972192
$ws.Range("Q19").Value = 571329.8719688904
$ws.Range("R19").Value = 5141968.847720014
$ws.Range("S19").Value = 0.2022271330569857
$ws.Range("T19").Value = 0.2022271330569858

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 2293.273345666667
$ws.Range("H20").Value = 6879.820037
$ws.Range("I20").Value = 0.5217768014597114
$ws.Range("J20").Value = 0.5217768014597114
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 186.9310863333334
$ws.Range("N20").Value = 560.793259
$ws.Range("O20").Value = 0.2908071500393791
$ws.Range("P20").Value = 0.2908071500393791
$ws.Range("Q20").Value = 428684.077764748
$ws.Range("R20").Value = 3858156.699882731
$ws.Range("S20").Value = 0.1517364245891616
$ws.Range("T20").Value = 0.1517364245891616

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 2293.273345666667
$ws.Range("H21").Value = 6879.820037
$ws.Range("I21").Value = 0.5217768014597114
$ws.Range("J21").Value = 0.5217768014597114
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 17.46044133333334
$ws.Range("N21").Value = 52.38132400000001
$ws.Range("O21").Value = 0.02716306464684043
$ws.Range("P21").Value = 0.02716306464684043
$ws.Range("Q21").Value = 40041.5647133099
$ws.Range("R21").Value = 360374.082419789
$ws.Range("S21").Value = 0.01417305698927176
$ws.Range("T21").Value = 0.01417305698927177

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 6.448603666666666
$ws.Range("H22").Value = 19.345811
$ws.Range("I22").Value = 0.00146721794043115
$ws.Range("J22").Value = 0.00146721794043115
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 6.744736333333333
$ws.Range("N22").Value = 20.234209
$ws.Range("O22").Value = 0.01049273071342527
$ws.Range("P22").Value = 0.01049273071342527
$ws.Range("Q22").Value = 43.49413144983322
$ws.Range("R22").Value = 391.447183048499
$ws.Range("S22").Value = 0.0000153951227468505
$ws.Range("T22").Value = 0.0000153951227468505

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 6.448603666666666
$ws.Range("H23").Value = 19.345811
$ws.Range("I23").Value = 0.00146721794043115
$ws.Range("J23").Value = 0.00146721794043115
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 182.5316823333334
$ws.Range("N23").Value = 547.595047
$ws.Range("O23").Value = 0.283963033503136
$ws.Range("P23").Value = 0.2839630335031361
$ws.Range("Q23").Value = 1177.074475977569
$ws.Range("R23").Value = 10593.67028379812
$ws.Range("S23").Value = 0.0004166356571750529
$ws.Range("T23").Value = 0.000416635657175053

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 6.448603666666666
$ws.Range("H24").Value = 19.345811
$ws.Range("I24").Value = 0.00146721794043115
$ws.Range("J24").Value = 0.00146721794043115
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 249.1329143333333
$ws.Range("N24").Value = 747.398743
$ws.Range("O24").Value = 0.3875740210972192
$ws.Range("P24").Value = 0.3875740210972192
$ws.Range("Q24").Value = 1606.559424857286
$ws.Range("R24").Value = 14459.03482371557
$ws.Range("S24").Value = 0.000568655556998881
$ws.Range("T24").Value = 0.0005686555569988811

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 6.448603666666666
$ws.Range("H25").Value = 19.345811
$ws.Range("I25").Value = 0.00146721794043115
$ws.Range("J25").Value = 0.00146721794043115
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 186.9310863333334
$ws.Range("N25").Value = 560.793259
$ws.Range("O25").Value = 0.2908071500393791
$ws.Range("P25").Value = 0.2908071500393791
$ws.Range("Q25").Value = 1205.444488743117
$ws.Range("R25").Value = 10849.00039868805
$ws.Range("S25").Value = 0.0004266774677434302
$ws.Range("T25").Value = 0.0004266774677434302

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 6.448603666666666
$ws.Range("H26").Value = 19.345811
$ws.Range("I26").Value = 0.00146721794043115
$ws.Range("J26").Value = 0.00146721794043115
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 17.46044133333334
$ws.Range("N26").Value = 52.38132400000001
$ws.Range("O26").Value = 0.02716306464684043
$ws.Range("P26").Value = 0.02716306464684043
$ws.Range("Q26").Value = 112.5954660037516
$ws.Range("R26").Value = 1013.359194033764
$ws.Range("S26").Value = 0.0000398541357669354
$ws.Range("T26").Value = 0.0000398541357669354

